# ManageClaims.xlsx - "Updation in custom order for both QA and STG"
#
# The workbook keeps a short rolling window of sample "custom order" test
# data on the ShipmentInformation / Input / ClaimDetail sheets (row 2 is
# always the most-recently-added order). This change rotates in a brand
# new order (PickUp90 / DropOff771, claim 59071590 dated 05-26-2022,
# tracking 999U684759, FCPBID1034621) in place of the previous one.

$wb = $excel.ActiveWorkbook

# ---- ShipmentInformation: pick-up / drop-off company name for the order
$wsShip = $wb.Worksheets.Item("ShipmentInformation")
$wsShip.Range("C2").Value = "PickUp90"
$wsShip.Range("K2").Value = "DropOff771"

# ---- Input: claim/order identifying fields for the order
$wsInput = $wb.Worksheets.Item("Input")
# Leading "'" forces these to stay plain text (not get reinterpreted as a
# date / number by Excel) just like the existing values in the column.
$wsInput.Range("B2").Formula = "'05-26-2022"
$wsInput.Range("T2").Formula = "'59071590"
$wsInput.Range("W2").Value = "999U684759"
$wsInput.Range("X2").Value = "FCPBID1034621"

# ---- ClaimDetail: widen column C (ClaimStatus) to fit the refreshed data
$wsClaim = $wb.Worksheets.Item("ClaimDetail")
$wsClaim.Columns.Item(3).ColumnWidth = 27.6
